# Update cryptocurrency price/volume data per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.721.45"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "1.543.72"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").Value = "'205.90"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("E6").Value = "  -1.39%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'21.41"
$ws.Range("E8").Value = "  -2.97%  "

$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("D10").Value = "'0.0580"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").Value = "1.766.15"
$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("D13").Value = "1.548.01"
$ws.Range("E13").Value = "  -1.70%  "

$ws.Range("D14").Value = "'3.67"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").Value = "'0.510"
$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").Value = "26.702.15"
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "'61.15"
$ws.Range("E17").Value = "  -1.01%  "

$ws.Range("D18").Value = "'212.39"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("D20").Value = "'7.22"
$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").Value = "'4.06"
$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").Value = "'8.94"
$ws.Range("E23").Value = "  -4.77%  "

$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").Value = "'152.15"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'14.86"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'6.47"
$ws.Range("E27").Value = "  -2.81%  "

$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("E29").Value = "  -1.00%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.10"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0459"
$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("D33").Value = "1.346.19"
$ws.Range("E33").Value = "  -3.02%  "

$ws.Range("D34").Value = "'2.91"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("E35").Value = "  -2.60%  "

$ws.Range("D36").Value = "'2.28"
$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("D37").Value = "'0.933"
$ws.Range("E37").Value = "  -1.23%  "

$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("E39").Value = "  +2.41%  "

$ws.Range("D40").Value = "'5.75"
$ws.Range("E40").Value = "  +5.99%  "

$ws.Range("D41").Value = "'0.799"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").Value = "'0.994"
$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.74"
$ws.Range("E44").Value = "  -3.51%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'62.51"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").Value = "1.679.26"
$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("E47").Value = "  -4.24%  "

$ws.Range("D48").Value = "'85.81"
$ws.Range("E48").Value = "  +0.66%  "

$ws.Range("E49").Value = "  +1.95%  "

$ws.Range("D50").Value = "0.0₇0977"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("E51").Value = "  +0.56%  "
